$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-10 (columns D, L, M, N, O, P, Q, R, S, T) are cyclically
# permuted. New row N gets the old values that used to live in row Map[N].
# Mapping expressed as newRow -> oldRow:
#   2 <- 7, 3 <- 2, 4 <- 10, 5 <- 8, 6 <- 9, 7 <- 6, 8 <- 3, 9 <- 4, 10 <- 5

$rows = @{
    2  = @{ D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 750;  T = 18 }
    3  = @{ D = 45092; L = "Primera"; M = 35;  N = 18000; O = 19000; P = 18571; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";     S = 1032; T = 18 }
    4  = @{ D = 44330; L = "Primera"; M = 60;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 861;  T = 18 }
    5  = @{ D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins";  S = 1109; T = 16 }
    6  = @{ D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins";  S = 1000; T = 16 }
    7  = @{ D = 45096; L = "Primera"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";     S = 1000; T = 18 }
    8  = @{ D = 45086; L = "Primera"; M = 30;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";     S = 1000; T = 18 }
    9  = @{ D = 44698; L = "Primera"; M = 120; N = 16000; O = 17000; P = 16500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins";  S = 917;  T = 18 }
    10 = @{ D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Región de O'Higgins"; S = 1042; T = 12 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
